# Add a new "2022-Q4" sheet of fund-holding detail (copied structure/style
# from the existing "2022-Q3" sheet, placed right before it in tab order),
# and insert the matching summary row at the top of the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Detail sheet: duplicate "2022-Q3" (same header/style layout as every
#    other quarterly sheet), drop the copy in before it, rename, and
#    overwrite its single data row with the 2022-Q4 figures.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

$q4.Cells.Item(2, 1).Value = 0

# Columns B-G hold numeric-looking figures that the source data keeps as
# plain text. Force text storage via a "@" number format, then reset the
# cell's style back to Normal so no stray formatting sticks around.
$textCols = 2, 3, 4, 5, 6, 7
$textVals = "260115", "景顺长城中小盘混合", "1.06", "91.71", "4.98", "0.0528"
for ($i = 0; $i -lt $textCols.Length; $i++) {
    $c = $q4.Cells.Item(2, $textCols[$i])
    $c.NumberFormat = "@"
    $c.Value = $textVals[$i]
    $c.Style = "Normal"
}

$q4.Cells.Item(2, 8).Value = 5

# ---------------------------------------------------------------------
# 2) Summary sheet "总计": insert a new row 2 (pushing the existing
#    quarters down by one), copy the formatting from the row that used
#    to be row 2 (now row 3) so the inserted row matches the table's
#    look, then fill in the 2022-Q4 totals.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 1
$total.Cells.Item(2, 4).Value = 0.05

# Column A is a plain 0-based row counter; re-sequence it for every row
# that got pushed down by the insert above.
for ($r = 3; $r -le 9; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}
